$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column C for rows 2-7
# from 2023-10-05 (45204) to 2023-10-08 (45207)
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = "2023-10-08"
}
